$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column cells are stored as text so values like "1.003" are preserved exactly
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "24.264.44"
$ws.Range("E2").Value = "  +14.13%  "
$ws.Range("D3").Value = "1.673.73"
$ws.Range("E3").Value = "  +8.21%  "
$ws.Range("D4").Value = "1.003"
$ws.Range("E4").Value = "  -0.61%  "
$ws.Range("D5").Value = "308.38"
$ws.Range("E5").Value = "  +9.24%  "
$ws.Range("D6").Value = "0.9984"
$ws.Range("E6").Value = "  +3.08%  "
$ws.Range("D7").Value = "0.3730"
$ws.Range("E7").Value = "  +2.87%  "
$ws.Range("D8").Value = "0.3426"
$ws.Range("E8").Value = "  +6.59%  "
$ws.Range("D9").Value = "47.46"
$ws.Range("E9").Value = "  +15.39%  "
$ws.Range("D10").Value = "1.183"
$ws.Range("E10").Value = "  +6.37%  "
$ws.Range("D11").Value = "0.07295"
$ws.Range("E11").Value = "  +5.27%  "
$ws.Range("D12").Value = "0.9989"
$ws.Range("E12").Value = "  -0.62%  "
$ws.Range("D13").Value = "20.43"
$ws.Range("E13").Value = "  +8.11%  "
$ws.Range("D14").Value = "6.098"
$ws.Range("E14").Value = "  +6.56%  "
$ws.Range("D15").Value = "6.761"
$ws.Range("E15").Value = "  +5.36%  "
$ws.Range("D16").Value = "1.671.54"
$ws.Range("E16").Value = "  +8.31%  "
$ws.Range("D17").Value = "0.00001108"
$ws.Range("E17").Value = "  +5.43%  "
$ws.Range("D18").Value = "0.9983"
$ws.Range("E18").Value = "  +3.26%  "
$ws.Range("D19").Value = "0.06721"
$ws.Range("E19").Value = "  +9.38%  "
$ws.Range("D20").Value = "81.61"
$ws.Range("E20").Value = "  +11.67%  "
$ws.Range("D21").Value = "16.44"
$ws.Range("E21").Value = "  +7.90%  "
$ws.Range("D22").Value = "6.136"
$ws.Range("E22").Value = "  +6.72%  "
$ws.Range("D23").Value = "12.00"
$ws.Range("E23").Value = "  +5.49%  "
$ws.Range("D24").Value = "24.213.22"
$ws.Range("E24").Value = "  +13.00%  "
$ws.Range("D25").Value = "2.411"
$ws.Range("E25").Value = "  +3.87%  "
$ws.Range("D26").Value = "3.357"
$ws.Range("E26").Value = "  -9.54%  "
$ws.Range("D27").Value = "2.658"
$ws.Range("E27").Value = "  +16.55%  "
$ws.Range("D28").Value = "151.81"
$ws.Range("E28").Value = "  +2.73%  "
$ws.Range("D29").Value = "19.50"
$ws.Range("E29").Value = "  +8.91%  "
$ws.Range("D30").Value = "1.855.99"
$ws.Range("E30").Value = "  +8.28%  "
$ws.Range("D31").Value = "127.17"
$ws.Range("E31").Value = "  +6.74%  "
$ws.Range("D32").Value = "6.404"
$ws.Range("E32").Value = "  +21.67%  "
$ws.Range("D33").Value = "4.139"
$ws.Range("E33").Value = "  +1.77%  "
$ws.Range("D34").Value = "0.9914"
$ws.Range("E34").Value = "  +14.05%  "
$ws.Range("D35").Value = "1.765"
$ws.Range("E35").Value = "  +15.82%  "
$ws.Range("D36").Value = "0.08455"
$ws.Range("E36").Value = "  +4.82%  "
$ws.Range("D37").Value = "12.59"
$ws.Range("E37").Value = "  +17.34%  "
$ws.Range("D38").Value = "0.06444"
$ws.Range("E38").Value = "  +9.69%  "
$ws.Range("D39").Value = "5.368"
$ws.Range("E39").Value = "  +7.66%  "
$ws.Range("D40").Value = "8.807"
$ws.Range("E40").Value = "  +11.87%  "
$ws.Range("D41").Value = "0.02346"
$ws.Range("E41").Value = "  +10.65%  "
$ws.Range("D43").Value = "0.2114"
$ws.Range("E43").Value = "  +9.39%  "
$ws.Range("D44").Value = "0.6186"
$ws.Range("E44").Value = "  +12.23%  "
$ws.Range("D45").Value = "0.9979"
$ws.Range("E45").Value = "  +3.21%  "
$ws.Range("D46").Value = "13.30"
$ws.Range("E46").Value = "  +5.08%  "
$ws.Range("D47").Value = "3.799"
$ws.Range("E47").Value = "  +6.24%  "
$ws.Range("D48").Value = "0.5952"
$ws.Range("E48").Value = "  +8.43%  "
$ws.Range("D49").Value = "127.18"
$ws.Range("E49").Value = "  +3.77%  "
$ws.Range("D50").Value = "2.028"
$ws.Range("E50").Value = "  +7.66%  "
$ws.Range("D51").Value = "0.07163"
$ws.Range("E51").Value = "  +8.10%  "
